# 4.0.3 model and data
# Expand the single "BVTQaZ" and "VTQaZ" transportation input-file rows on the
# "Boolean" sheet into six mode-specific rows each (LDVs, HDVs, aircraft,
# rail, ships, motorbikes), matching the new upstream CSV file layout.

$wb = $excel.ActiveWorkbook

$wsBoolean = $wb.Worksheets.Item("Boolean")

# --- Expand "trans/BVTQaZ/BVTQaZ.csv" (row 17) into 6 rows ---
# Insert 5 new rows below row 17 so rows 17-22 are available.
$wsBoolean.Range("A18:A22").EntireRow.Insert()

$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the first insertion, the old row 21 ("trans/VTQaZ/VTQaZ.csv") is
# now at row 26 (21 + 5).
$wsBoolean.Range("A27:A31").EntireRow.Insert()

$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# --- Trailing blank formatted rows (33-38) ---
for ($r = 33; $r -le 38; $r++) {
  $wsBoolean.Rows.Item($r).Font.Name = "Calibri"
}

# --- View / selection bookkeeping to mirror the saved workbook state ---
$wsInteger = $wb.Worksheets.Item("Integer")
$wsInteger.Activate()
$wsInteger.Range("A13").Select()

$wsBoolean.Activate()
$wsBoolean.Range("A32").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

Write-Host "done"
